$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4201.3076
$ws.Range("J38").Value = 6177.5
$ws.Range("L38").Value = 18532.5
$ws.Range("N38").Value = -19276.5

$ws.Range("H51").Value = 6002
$ws.Range("J51").Value = 6002
$ws.Range("L51").Value = 6002
$ws.Range("N51").Value = -6970

$ws.Range("H74").Value = 2299
$ws.Range("I74").Value = 999
$ws.Range("J74").Value = 3599
$ws.Range("K74").Value = 999
$ws.Range("L74").Value = 3599
$ws.Range("M74").Value = -63
$ws.Range("N74").Value = -5471

$ws.Range("H77").Value = 2299
$ws.Range("I77").Value = 999
$ws.Range("J77").Value = 3599
$ws.Range("K77").Value = 4995
$ws.Range("L77").Value = 17995
$ws.Range("M77").Value = -315
$ws.Range("N77").Value = -27355

$ws.Range("H88").Value = 807.0909
$ws.Range("I88").Value = 789.8
$ws.Range("K88").Value = 789.8
$ws.Range("M88").Value = -383.8

$ws.Range("H91").Value = 807.0909
$ws.Range("I91").Value = 789.8
$ws.Range("K91").Value = 789.8
$ws.Range("M91").Value = 614.2

$ws.Range("H99").Value = 3977
$ws.Range("I99").Value = 2707
$ws.Range("K99").Value = 8121
$ws.Range("M99").Value = -6623

$ws.Range("H100").Value = 2390
$ws.Range("I100").Value = 2390
$ws.Range("K100").Value = 2390
$ws.Range("M100").Value = -1849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 30000
$ws.Range("J23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("N23").Value = -30518

$ws.Range("H45").Value = 2472.3333
$ws.Range("I45").Value = 1493.8334
$ws.Range("K45").Value = 1493.8334
$ws.Range("M45").Value = -1116.8334

$ws.Range("H122").Value = 1342
$ws.Range("I122").Value = 1006
$ws.Range("K122").Value = 3018
$ws.Range("M122").Value = -568

$ws.Range("H132").Value = 3886.3684
$ws.Range("I132").Value = 3797.5
$ws.Range("J132").Value = 4360.3335
$ws.Range("K132").Value = 11392.5
$ws.Range("L132").Value = 13081.0005
$ws.Range("M132").Value = -8862.5
$ws.Range("N132").Value = -18141.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5472.2144
$ws.Range("I107").Value = 4734.6665
$ws.Range("K107").Value = 4734.6665
$ws.Range("M107").Value = -2814.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4762.5
$ws.Range("I62").Value = 4762.5
$ws.Range("K62").Value = 4762.5
$ws.Range("M62").Value = -4138.5

$ws.Range("H65").Value = 4762.5
$ws.Range("I65").Value = 4762.5
$ws.Range("K65").Value = 23812.5
$ws.Range("M65").Value = -20692.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 8746
$ws.Range("I9").Value = 7500
$ws.Range("J9").Value = 9992
$ws.Range("K9").Value = 22500
$ws.Range("L9").Value = 29976
$ws.Range("M9").Value = -22276
$ws.Range("N9").Value = -30424

$ws.Range("H34").Value = 721.4666999999999
$ws.Range("J34").Value = 1037.3
$ws.Range("L34").Value = 3111.9
$ws.Range("N34").Value = -3279.9

$ws.Range("H46").Value = 200
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

$ws.Range("H57").Value = 1362.6842
$ws.Range("I57").Value = 1298.875
$ws.Range("J57").Value = 1409.091
$ws.Range("K57").Value = 3896.625
$ws.Range("L57").Value = 4227.272999999999
$ws.Range("M57").Value = -3337.625
$ws.Range("N57").Value = -5345.272999999999

$ws.Range("H58").Value = 2749.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2749.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8248.5
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -8504.5

$ws.Range("H81").Value = 1119.1428
$ws.Range("I81").Value = 968
$ws.Range("K81").Value = 2904
$ws.Range("M81").Value = -1781

$ws.Range("H84").Value = 1119.1428
$ws.Range("I84").Value = 968
$ws.Range("K84").Value = 8712
$ws.Range("M84").Value = -3096

$ws.Range("H86").Value = 421.1111
$ws.Range("I86").Value = 424.125
$ws.Range("J86").Value = 397
$ws.Range("K86").Value = 1272.375
$ws.Range("L86").Value = 1191
$ws.Range("M86").Value = -86.375
$ws.Range("N86").Value = -3563

$ws.Range("H89").Value = 421.1111
$ws.Range("I89").Value = 424.125
$ws.Range("J89").Value = 397
$ws.Range("K89").Value = 3817.125
$ws.Range("L89").Value = 3573
$ws.Range("M89").Value = 2110.875
$ws.Range("N89").Value = -15429

$ws.Range("H109").Value = 2486.125
$ws.Range("I109").Value = 1786.3334
$ws.Range("J109").Value = 2906
$ws.Range("K109").Value = 5359.0002
$ws.Range("L109").Value = 8718
$ws.Range("M109").Value = -4319.0002
$ws.Range("N109").Value = -10798

$ws.Range("H112").Value = 675.5
$ws.Range("I112").Value = 643
$ws.Range("K112").Value = 1929
$ws.Range("M112").Value = -821

$ws.Range("H121").Value = 1180
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 1766.6666
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 5299.9998
$ws.Range("M121").Value = 410
$ws.Range("N121").Value = -7919.9998

$ws.Range("H131").Value = 1635.1
$ws.Range("I131").Value = 873.6667
$ws.Range("K131").Value = 2621.0001
$ws.Range("M131").Value = 2418.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7500.231
$ws.Range("I113").Value = 6660.625
$ws.Range("J113").Value = 8843.6
$ws.Range("K113").Value = 6660.625
$ws.Range("L113").Value = 8843.6
$ws.Range("M113").Value = -4490.625
$ws.Range("N113").Value = -13183.6

$ws.Range("H123").Value = 983333.3
$ws.Range("J123").Value = 950000
$ws.Range("L123").Value = 950000
$ws.Range("N123").Value = -954900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null

$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = $null

$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = $null

$ws.Range("H40").Value = 4999
$ws.Range("I40").Value = 4999
$ws.Range("K40").Value = 4999
$ws.Range("M40").Value = -4863

$ws.Range("H46").Value = 4817.25
$ws.Range("I46").Value = 2137
$ws.Range("J46").Value = 7497.5
$ws.Range("K46").Value = 2137
$ws.Range("L46").Value = 7497.5
$ws.Range("M46").Value = -1949
$ws.Range("N46").Value = -7873.5

$ws.Range("H82").Value = 4040.9092
$ws.Range("I82").Value = 625
$ws.Range("J82").Value = 4800
$ws.Range("K82").Value = 625
$ws.Range("L82").Value = 4800
$ws.Range("M82").Value = -264
$ws.Range("N82").Value = -5522

$ws.Range("H85").Value = 4040.9092
$ws.Range("I85").Value = 625
$ws.Range("J85").Value = 4800
$ws.Range("K85").Value = 625
$ws.Range("L85").Value = 4800
$ws.Range("M85").Value = 623
$ws.Range("N85").Value = -7296

$ws.Range("H125").Value = 35000
$ws.Range("J125").Value = 35000
$ws.Range("L125").Value = 35000
$ws.Range("N125").Value = -44840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 75000
$ws.Range("I4").Value = 100000
$ws.Range("J4").Value = 50000
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 50000
$ws.Range("M4").Value = -99887
$ws.Range("N4").Value = -50226

$ws.Range("H107").Value = 730
$ws.Range("I107").Value = 583.6667
$ws.Range("K107").Value = 1751.0001
$ws.Range("M107").Value = 168.9999

$ws.Range("H122").Value = 3074.5
$ws.Range("I122").Value = 1661.75
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 4985.25
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -2535.25
$ws.Range("N122").Value = -22600

$ws.Range("H132").Value = 1775.7142
$ws.Range("I132").Value = 1714.5
$ws.Range("K132").Value = 5143.5
$ws.Range("M132").Value = -2613.5
